$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'26.879.60"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -2.00%  '
$ws.Range('D3').Value = "'1.815.83"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.10%  '
$ws.Range('E4').Value = '  -0.45%  '
$ws.Range('E5').Value = '  -0.36%  '
$ws.Range('D6').Value = "'308.56"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.84%  '
$ws.Range('E7').Value = '  -2.53%  '
$ws.Range('D8').Value = "'0.3652"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.07%  '
$ws.Range('D9').Value = "'0.07228"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.15%  '
$ws.Range('D10').Value = "'0.8576"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.26%  '
$ws.Range('D11').Value = "'19.71"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.75%  '
$ws.Range('D12').Value = "'0.07535"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.65%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = "'5.318"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.50%  '
$ws.Range('B14').Value = 'Litecoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D14').Value = "'91.75"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.46%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = "'1.719.75"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -9.18%  '
$ws.Range('D16').Value = "'6.470"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.71%  '
$ws.Range('E17').Value = '  -0.27%  '
$ws.Range('D18').Value = "'0.000008620"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.32%  '
$ws.Range('D20').Value = "'14.43"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.51%  '
$ws.Range('D21').Value = "'26.672.54"
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').Value = "'5.126"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.61%  '
$ws.Range('E23').Value = '  -1.74%  '
$ws.Range('D24').Value = "'1.951.98"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -7.15%  '
$ws.Range('D25').Value = "'151.69"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.38%  '
$ws.Range('D26').Value = "'1.842"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.35%  '
$ws.Range('E27').Value = '  -2.91%  '
$ws.Range('D28').Value = "'2.076"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.00%  '
$ws.Range('D30').Value = "'114.93"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.26%  '
$ws.Range('D31').Value = "'0.08842"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.71%  '
$ws.Range('D32').Value = "'2.958"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.50%  '
$ws.Range('D33').Value = "'4.405"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.26%  '
$ws.Range('D34').Value = "'1.127"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.28%  '
$ws.Range('D35').Value = "'0.7137"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -5.86%  '
$ws.Range('E36').Value = '  -2.51%  '
$ws.Range('D37').Value = "'0.05239"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.06%  '
$ws.Range('D38').Value = "'2.419"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.22%  '
$ws.Range('D39').Value = "'0.01914"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.11%  '
$ws.Range('D40').Value = "'2.917"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.13%  '
$ws.Range('D41').Value = "'7.138"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.65%  '
$ws.Range('D42').Value = "'0.5139"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.65%  '
$ws.Range('D43').Value = "'0.1621"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.25%  '
$ws.Range('E44').Value = '  -4.25%  '
$ws.Range('D45').Value = "'0.4797"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.57%  '
$ws.Range('E46').Value = '  -0.43%  '
$ws.Range('B47').Value = 'Quant'
$ws.Range('C47').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D47').Value = "'103.01"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.85%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = "'10.05"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.68%  '
$ws.Range('D49').Value = "'0.06282"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.39%  '
$ws.Range('D50').Value = "'1.615"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.71%  '
$ws.Range('D51').Value = "'63.82"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.06%  '
